# Weekly update: insert a new price-report row (date 2023-02-03 / serial 44960)
# above the existing "Provincia de Diguillín" block, pushing the older rows
# down by one (row 16 -> 22 become rows 17 -> 22, plus a brand-new row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 16; everything currently at 16..21 shifts to 17..22.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with this week's record.
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 44960
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100101
$ws.Range("H16").Value = "Berries"
$ws.Range("I16").Value = 100101001
$ws.Range("J16").Value = "Arándano (blue)"
$ws.Range("K16").Value = "Sin especificar"
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 2500
$ws.Range("O16").Value = 2500
$ws.Range("P16").Value = 2500
$ws.Range("Q16").Value = "$/bandeja 2 kilos"
$ws.Range("R16").Value = "Provincia de Diguillín"
$ws.Range("S16").Value = 1250
$ws.Range("T16").Value = 2
